$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at 25 (pushes the pik-potsdam... "solve()" block down by
#    one row) and fill it with the first newly-added URL.
# ---------------------------------------------------------------------------
$ws.Rows("25:25").Insert()
$ws.Cells.Item(25, 1).Value = "https://www.umweltbundesamt.de/sites/default/files/medien/1410/publikationen/2020-12-21_methodenkonvention_3_1_kostensaetze.pdf"

# ---------------------------------------------------------------------------
# 2) Write the second newly-added URL straight into row 34. Rows 35-55 are
#    empty in the sheet (sparse), so this does not need an Insert and the
#    existing "solve()" row (56 -> 57 because of step 1 only) is unaffected.
# ---------------------------------------------------------------------------
$ws.Cells.Item(34, 1).Value = "https://www.waermepumpe.de/uploads/tx_bcpageflip/BWP_Branchenstudie2021_Update.pdf"

# ---------------------------------------------------------------------------
# 3) The row-insert in step 1 does not auto-shift the worksheet Hyperlinks
#    collection, so the existing hyperlinks need to be re-pointed at their
#    new cells. This engine's Hyperlinks collection only supports whole-
#    collection Delete + Add, so: snapshot each hyperlinked cell's
#    formatting, delete every hyperlink, re-add them (with the one cell that
#    shifted, A28 -> A29, corrected), then restore formatting.
# ---------------------------------------------------------------------------
$links = @(
    @{ Ref = "A6";  Target = "https://www.heizungsdiscount24.de/waermepumpen/vaillant-versotherm-plus-vwl-775-luft-wasser-waermepumpe.html?gclid=CjwKCAjwnZaVBhA6EiwAVVyv9MQZvSx9QQuF56cGi9y1Cq8h1lNVzaH_q0FYiCaP7LpHmW8Vs_3EeBoCkU4QAvD_BwE&cq_cmp=13242830342&cq_plt=gp&cq_src=google_ads&cq_net=u"; SubAddress = ""; Display = "https://www.heizungsdiscount24.de/waermepumpen/vaillant-versotherm-plus-vwl-775-luft-wasser-waermepumpe.html?gclid=CjwKCAjwnZaVBhA6EiwAVVyv9MQZvSx9QQuF56cGi9y1Cq8h1lNVzaH_q0FYiCaP7LpHmW8Vs_3EeBoCkU4QAvD_BwE&cq_cmp=13242830342&cq_plt=gp&cq_src=google_ads&cq_net=u" },
    @{ Ref = "A11"; Target = "https://domotec.ch/wp-content/uploads/2022/06/1.1-pl-allgemein-06.2022-DE.pdf"; SubAddress = ""; Display = $null },
    @{ Ref = "A9";  Target = "https://docplayer.org/82079403-Preisliste-waermepumpen-systeme-der-cta-ag.html"; SubAddress = ""; Display = $null },
    @{ Ref = "A10"; Target = "https://shop.smuk.at/shop/USER_ARTIKEL_HANDLING_AUFRUF.php?Kategorie_ID=9389&Ziel_ID=12271890"; SubAddress = ""; Display = $null },
    @{ Ref = "A8";  Target = "https://www.heizungsdiscount24.de/waermepumpen/vaillant-flexotherm-exclusive-vwf-574-heizungswaermepumpe-solewasser.html?cq_src=google_ads&cq_net=u&cq_cmp=13242830342&cq_plt=gp&gclid=CjwKCAjwnZaVBhA6EiwAVVyv9LDV4ncrTDuayjy2mZ2XWxvqs-T0jg902k_jxM-pgcEy8--TXt17SRoCTbwQAvD_BwE"; SubAddress = ""; Display = "https://www.heizungsdiscount24.de/waermepumpen/vaillant-flexotherm-exclusive-vwf-574-heizungswaermepumpe-solewasser.html?cq_src=google_ads&cq_net=u&cq_cmp=13242830342&cq_plt=gp&gclid=CjwKCAjwnZaVBhA6EiwAVVyv9LDV4ncrTDuayjy2mZ2XWxvqs-T0jg902k_jxM-pgcEy8--TXt17SRoCTbwQAvD_BwE" },
    @{ Ref = "A7";  Target = "https://domotec.ch/wp-content/uploads/2022/06/1.1-pl-allgemein-06.2022-DE.pdfhttps:/domotec.ch/wp-content/uploads/2022/06/1.1-pl-allgemein-06.2022-DE.pdf"; SubAddress = ""; Display = $null },
    @{ Ref = "A5";  Target = "https://heizung-billiger.de/69503-stiebel-eltron-luft-wasser-warmepumpe-wpl-09-ikcs-classic-stiebel-236377-4017212363775.html?hbdc=DE&utm_source=guenstiger&utm_medium=cpc&utm_campaign=guenstiger-de"; SubAddress = ""; Display = $null },
    @{ Ref = "A12"; Target = "https://www.preis.de/katalog/Viessmann-Waermepumpen/14892.html"; SubAddress = ""; Display = $null },
    @{ Ref = "A18"; Target = "https://www.globalpetrolprices.com/Germany/natural_gas_prices/"; SubAddress = ""; Display = $null },
    @{ Ref = "A22"; Target = "https://www.umweltbundesamt.de/daten/umwelt-wirtschaft/gesellschaftliche-kosten-von-umweltbelastungen"; SubAddress = "klimakosten-von-treibhausgas-emissionen"; Display = $null },
    @{ Ref = "A29"; Target = "https://www.viessmann.de/content/dam/vi-brands/DE/PDF/Planungshandbuch/ph-waermepumpen.pdf/_jcr_content/renditions/original.media_file.download_attachment.file/ph-waermepumpen.pdf"; SubAddress = ""; Display = $null },
    @{ Ref = "A20"; Target = "https://www.eon.de/de/pk/strom/stromanbieter/guenstiger-stromanbieter.html"; SubAddress = ""; Display = $null },
    @{ Ref = "A21"; Target = "https://www.volker-quaschning.de/datserv/CO2-spez/index_e.php"; SubAddress = ""; Display = $null }
)

# Snapshot current formatting of every cell that will carry a hyperlink into
# a staging area far below the data, so it can be restored after the
# Hyperlinks collection is rebuilt (Hyperlinks.Add() re-styles the cell).
$stageRow = 1000
foreach ($it in $links) {
    $src = $ws.Range($it.Ref)
    $stageCell = $ws.Cells.Item($stageRow, 1)
    $src.Copy()
    $stageCell.PasteSpecial(-4122)
    $stageRow = $stageRow + 1
}

$ws.Hyperlinks.Delete()

$stageRow = 1000
foreach ($it in $links) {
    $dst = $ws.Range($it.Ref)
    if ($it.Display) {
        $ws.Hyperlinks.Add($dst, $it.Target, $it.SubAddress, "", $it.Display)
    } else {
        $ws.Hyperlinks.Add($dst, $it.Target, $it.SubAddress)
    }
    $stageCell = $ws.Cells.Item($stageRow, 1)
    $stageCell.Copy()
    $dst.PasteSpecial(-4122)
    $stageRow = $stageRow + 1
}

# Clean up the staging area.
$ws.Range("A1000:A1012").Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Update the view: scrolled so row 19 is at the top, with A34 selected.
# ---------------------------------------------------------------------------
$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 19
